$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A (col 1) from 17 to 22 "characters" as stored in the
# worksheet XML. The Excel object model's ColumnWidth setter bakes in a
# fixed +5/6 character padding before it lands in the saved <col width=.../>
# attribute (e.g. ColumnWidth=17 round-trips to width="17.8333..."), so we
# back that padding out here to land on the exact target width="22".
$ws.Columns.Item(1).ColumnWidth = 22 - (5/6)

# Rows 2 and 3 effectively swap their B-F contents (the process record that
# used to be "processo_2" moves up to row 2, and row 3 is relabeled as the
# new "processo_word" record instead of the old "copy_processo_1"). The OAB
# column (E) becomes a real number in both rows, matching the numeric type
# already used by E2 before the edit.

$ws.Range("A2").Value = "processo_2"
$ws.Range("B2").Value = "4835245-15.2024.8.01.2832"
$ws.Range("C2").Value = "Nome Aleatório 2"
$ws.Range("D2").Value = "Advogado Exemplo"
$ws.Range("E2").Value = 12723
$ws.Range("F2").Value = "25/5/2024"

$ws.Range("A3").Value = "processo_word"
$ws.Range("B3").Value = "3781128-20.2024.8.01.8252"
$ws.Range("C3").Value = "Nome Aleatório 86"
$ws.Range("D3").Value = "Advogado Exemplo"
$ws.Range("E3").Value = 44432

# "12/5/2024" (day 12 / month 5) is also a valid US-style month/day date
# ("December 5"), so a plain .Value assignment gets silently reinterpreted
# as a date serial number by Excel's input parser. Route it through a
# text formula + paste-values round trip instead, which keeps it as a
# literal text string (like the original inlineStr cell) without
# disturbing the cell's existing style.
$ws.Range("F3").Formula = "=""12/5/2024"""
$ws.Range("F3").Copy()
$ws.Range("F3").PasteSpecial(-4163)

$wb.Save()
